$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("df2")

# Current layout (before):
#   A1:Nr.  B1:varName  C1:varLabel  D1:keyword1  E1:keyword2
# Target layout (after) - keyword columns re-sorted alphabetically,
# with a new "keyword3" column inserted before the (now last) keyword1 column:
#   A1:Nr.  B1:varName  C1:varLabel  D1:keyword2  E1:keyword3  F1:keyword1
# and the "x" value that used to live in the (old) D column (keyword1, row2)
# moves along with its header to column F, plus a new x value in F4.

$ws.Range("D1").Value = "keyword2"
$ws.Range("E1").Value = "keyword3"
$ws.Range("F1").Value = "keyword1"

$ws.Range("F4").Value = "x"

$ws.Range("F1").Select()
